$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1120.6957
$ws.Range("J2").Value = 2780
$ws.Range("L2").Value = 2780
$ws.Range("N2").Value = -3006
$ws.Range("H9").Value = 783.82355
$ws.Range("I9").Value = 541.9231
$ws.Range("K9").Value = 541.9231
$ws.Range("M9").Value = -372.9231
$ws.Range("H17").Value = 3630.3809
$ws.Range("J17").Value = 3711.9
$ws.Range("L17").Value = 11135.7
$ws.Range("N17").Value = -11471.7
$ws.Range("H18").Value = 836.2
$ws.Range("I18").Value = 836.2
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 836.2
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -552.2
$ws.Range("N18").Value = ""
$ws.Range("H86").Value = 4566.3076
$ws.Range("I86").Value = 2237.8
$ws.Range("J86").Value = 6021.625
$ws.Range("K86").Value = 2237.8
$ws.Range("L86").Value = 6021.625
$ws.Range("M86").Value = -1114.8
$ws.Range("N86").Value = -8267.625
$ws.Range("H89").Value = 4566.3076
$ws.Range("I89").Value = 2237.8
$ws.Range("J89").Value = 6021.625
$ws.Range("K89").Value = 11189
$ws.Range("L89").Value = 30108.125
$ws.Range("M89").Value = -5573
$ws.Range("N89").Value = -41340.125
$ws.Range("H95").Value = 24000
$ws.Range("J95").Value = 24000
$ws.Range("L95").Value = 24000
$ws.Range("N95").Value = -29492
$ws.Range("H127").Value = 5289.778
$ws.Range("I127").Value = 5968
$ws.Range("K127").Value = 17904
$ws.Range("M127").Value = -12944

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 162.11111
$ws.Range("I5").Value = 174.625
$ws.Range("J5").Value = 62
$ws.Range("K5").Value = 174.625
$ws.Range("L5").Value = 62
$ws.Range("M5").Value = -62.625
$ws.Range("N5").Value = -286
$ws.Range("H6").Value = 1500
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = ""
$ws.Range("H24").Value = 90000
$ws.Range("J24").Value = 90000
$ws.Range("L24").Value = 90000
$ws.Range("N24").Value = -90748
$ws.Range("H32").Value = 6561.2666
$ws.Range("I32").Value = 6886.077
$ws.Range("K32").Value = 6886.077
$ws.Range("M32").Value = -6599.077
$ws.Range("H45").Value = 5893.1333
$ws.Range("I45").Value = 6199.769
$ws.Range("K45").Value = 6199.769
$ws.Range("M45").Value = -5822.769
$ws.Range("H63").Value = 27901
$ws.Range("I63").Value = 32376.25
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 32376.25
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -31690.25
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 27901
$ws.Range("I66").Value = 32376.25
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 161881.25
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -158449.25
$ws.Range("N66").Value = -56864
$ws.Range("H100").Value = 90000
$ws.Range("J100").Value = 90000
$ws.Range("L100").Value = 90000
$ws.Range("N100").Value = -92164
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 162.11111
$ws.Range("I4").Value = 174.625
$ws.Range("J4").Value = 62
$ws.Range("K4").Value = 174.625
$ws.Range("L4").Value = 62
$ws.Range("M4").Value = -59.625
$ws.Range("N4").Value = -292
$ws.Range("H95").Value = 14031
$ws.Range("J95").Value = 14031
$ws.Range("L95").Value = 14031
$ws.Range("N95").Value = -19523

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 1000
$ws.Range("M39").Value = -609
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 1000
$ws.Range("M49").Value = -818
$ws.Range("H58").Value = 1793.4546
$ws.Range("I58").Value = 1712.1904
$ws.Range("K58").Value = 1712.1904
$ws.Range("M58").Value = -1509.1904
$ws.Range("H94").Value = 84980.21000000001
$ws.Range("J94").Value = 7046.143
$ws.Range("L94").Value = 7046.143
$ws.Range("N94").Value = -7948.143
$ws.Range("H99").Value = 3098
$ws.Range("I99").Value = 2299.75
$ws.Range("K99").Value = 2299.75
$ws.Range("M99").Value = -801.75
$ws.Range("H122").Value = 3918.8
$ws.Range("J122").Value = 3948
$ws.Range("L122").Value = 11844
$ws.Range("N122").Value = -16744
$ws.Range("H126").Value = 3098
$ws.Range("I126").Value = 2299.75
$ws.Range("K126").Value = 6899.25
$ws.Range("M126").Value = -4429.25
$ws.Range("H136").Value = 1793.4546
$ws.Range("I136").Value = 1712.1904
$ws.Range("K136").Value = 5136.5712
$ws.Range("M136").Value = -2586.5712

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -524
$ws.Range("H122").Value = 900
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -13000

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 31333.334
$ws.Range("J95").Value = 31333.334
$ws.Range("L95").Value = 31333.334
$ws.Range("N95").Value = -36825.334
$ws.Range("H98").Value = 20533.334
$ws.Range("J98").Value = 20533.334
$ws.Range("L98").Value = 20533.334
$ws.Range("N98").Value = -26523.334
$ws.Range("H122").Value = 3998.25
$ws.Range("I122").Value = 3998.5
$ws.Range("K122").Value = 11995.5
$ws.Range("M122").Value = -9545.5
$ws.Range("H126").Value = 4665.6665
$ws.Range("I126").Value = 3999.5
$ws.Range("J126").Value = 4998.75
$ws.Range("K126").Value = 11998.5
$ws.Range("L126").Value = 14996.25
$ws.Range("M126").Value = -9528.5
$ws.Range("N126").Value = -19936.25

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4223.2856
$ws.Range("I22").Value = 3795.8333
$ws.Range("K22").Value = 3795.8333
$ws.Range("M22").Value = -3500.8333
$ws.Range("H27").Value = 4223.2856
$ws.Range("I27").Value = 3795.8333
$ws.Range("K27").Value = 3795.8333
$ws.Range("M27").Value = -3688.8333
$ws.Range("H40").Value = 3297.8333
$ws.Range("I40").Value = 2398.8
$ws.Range("K40").Value = 2398.8
$ws.Range("M40").Value = -2262.8
$ws.Range("H46").Value = 3696.75
$ws.Range("H114").Value = 22000
$ws.Range("J114").Value = 22000
$ws.Range("L114").Value = 22000
$ws.Range("N114").Value = -30678
$ws.Range("H127").Value = 82497.5
$ws.Range("J127").Value = 82497.5
$ws.Range("L127").Value = 82497.5
$ws.Range("N127").Value = -92417.5

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 20572
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 20572
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 20572
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -22554
